$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-16 20:44:20"

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("H3").Value = "2016-08-16 20:44:15"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e64b201b12995117d5413a262ebbdd30c5e30973/e2e/bd1f1b19-458e-48fd-ab2f-511a5367d156.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a681bca58520842ce8b8af461ae244014679d763/e2e/bd1f1b19-458e-48fd-ab2f-511a5367d156.md."
$ws2.Columns.Item(16).ColumnWidth = 39.17

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("H3").Value = "2016-08-16 20:44:20"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e64b201b12995117d5413a262ebbdd30c5e30973/e2e/bd1f1b19-458e-48fd-ab2f-511a5367d156.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a681bca58520842ce8b8af461ae244014679d763/e2e/bd1f1b19-458e-48fd-ab2f-511a5367d156.md."
$ws3.Columns.Item(16).ColumnWidth = 39.17
